$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update the sample number text "E7420" -> "E7420L" in column G (G2:G39)
$ws.Range("G2:G39").Value = "E7420L"

# 2. Reset scroll position of the sheet view so the top-left visible cell is A1
#    (selection itself stays on G2:G39, only the scrolled/frozen top-left cell moves)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1

# 3. Replace the =FALSE() formulas in H2:H39 with a plain boolean FALSE value
$ws.Range("H2:H39").Value = $false
